$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.194.21"
$ws.Range("E2").Value = "  +1.89%  "

$ws.Range("D3").Value = "3.470.31"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'581.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'146.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("D7").Value = "3.469.22"
$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").Value = "'7.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("E12").Value = "  +4.52%  "

$ws.Range("D13").Value = "4.065.87"
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").Value = "'29.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("E15").Value = "  +2.30%  "

$ws.Range("D16").Value = "3.477.01"
$ws.Range("E16").Value = "  +1.48%  "

$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "63.250.55"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("E19").Value = "  +2.98%  "

$ws.Range("D20").Value = "'14.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.32%  "

$ws.Range("D21").Value = "'9.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("D22").Value = "'388.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "

$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("D24").Value = "'74.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "3.620.95"
$ws.Range("E26").Value = "  +1.66%  "

$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("E28").Value = "  -5.01%  "

$ws.Range("D29").Value = "'7.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.73%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "'8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -3.41%  "

$ws.Range("D35").Value = "'23.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").Value = "'5.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("E37").Value = "  +2.49%  "

$ws.Range("E38").Value = "  +8.16%  "

$ws.Range("D39").Value = "'31.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.38%  "

$ws.Range("D40").Value = "'169.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("D41").Value = "3.509.11"
$ws.Range("E41").Value = "  +1.43%  "

$ws.Range("D42").Value = "'0.0764"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").Value = "'0.798"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.28%  "

$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("D45").Value = "'42.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("E46").Value = "  +2.72%  "

$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("D48").Value = "2.604.29"
$ws.Range("E48").Value = "  +3.46%  "

$ws.Range("D49").Value = "'2.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.95%  "

$ws.Range("D50").Value = "'23.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").Value = "'6.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.69%  "
